{"js": "// The document contains a single 20x5 table of arithmetic problems\n// (\"a-b=c\" / \"a+b=c\" style strings), one per cell, each cell holding a\n// single paragraph with a single run. The edit replaces every one of the\n// 100 problem strings with a new one (same position in the table),\n// leaving all other content (the date paragraph, table/cell formatting,\n// fonts, sizes, alignment, etc.) untouched.\n\nconst newValues = [\n  [\"56-22=34\", \"98-49=49\", \"84-34=50\", \"75-29=46\", \"50+42=92\"],\n  [\"49-0=49\", \"94-59=35\", \"96-41=55\", \"20-10=10\", \"78-59=19\"],\n  [\"38+1=39\", \"74-50=24\", \"35+9=44\", \"84-81=3\", \"89-89=0\"],\n  [\"54-27=27\", \"34+44=78\", \"3+96=99\", \"7+4=11\", \"62-30=32\"],\n  [\"64+26=90\", \"76-73=3\", \"73-12=61\", \"5+35=40\", \"43-15=28\"],\n  [\"58-0=58\", \"24+48=72\", \"50+8=58\", \"48+38=86\", \"72-8=64\"],\n  [\"16+2=18\", \"39+33=72\", \"98-93=5\", \"51-23=28\", \"81+5=86\"],\n  [\"30+67=97\", \"55+32=87\", \"74-62=12\", \"65+30=95\", \"10+49=59\"],\n  [\"20+61=81\", \"30+8=38\", \"9+22=31\", \"53+24=77\", \"8+54=62\"],\n  [\"29+19=48\", \"83-75=8\", \"91-78=13\", \"14+14=28\", \"33+44=77\"],\n  [\"21+40=61\", \"21+24=45\", \"77+0=77\", \"91-43=48\", \"46-9=37\"],\n  [\"89-66=23\", \"77+0=77\", \"97-9=88\", \"19+25=44\", \"3+24=27\"],\n  [\"54-7=47\", \"58+8=66\", \"10+26=36\", \"84-50=34\", \"82-75=7\"],\n  [\"25+46=71\", \"37+26=63\", \"91-6=85\", \"25+44=69\", \"73-29=44\"],\n  [\"12+65=77\", \"45+37=82\", \"54+45=99\", \"80-44=36\", \"18-12=6\"],\n  [\"63+36=99\", \"47+30=77\", \"34+39=73\", \"76-65=11\", \"88-75=13\"],\n  [\"61-29=32\", \"83+4=87\", \"4+11=15\", \"69+30=99\", \"51+15=66\"],\n  [\"19+60=79\", \"61-43=18\", \"93-93=0\", \"54+26=80\", \"64-37=27\"],\n  [\"2+19=21\", \"87-8=79\", \"9+61=70\", \"42+31=73\", \"80-55=25\"],\n  [\"99-6=93\", \"59+9=68\", \"34+59=93\", \"94-32=62\", \"78-4=74\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document contains a single 20x5 table of arithmetic problems\n# (\"a-b=c\" / \"a+b=c\" style strings), one per cell. The edit replaces every\n# one of the 100 problem strings with a new one (same row/column position),\n# leaving all other content (the date paragraph, table/cell formatting,\n# fonts, sizes, alignment, etc.) untouched. Setting Cell(r,c).Range.Text\n# replaces only the cell's text run while Word keeps the existing\n# paragraph/run formatting and the cell's end-of-cell marker.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$values = @(\n    @(\"56-22=34\", \"98-49=49\", \"84-34=50\", \"75-29=46\", \"50+42=92\"),\n    @(\"49-0=49\", \"94-59=35\", \"96-41=55\", \"20-10=10\", \"78-59=19\"),\n    @(\"38+1=39\", \"74-50=24\", \"35+9=44\", \"84-81=3\", \"89-89=0\"),\n    @(\"54-27=27\", \"34+44=78\", \"3+96=99\", \"7+4=11\", \"62-30=32\"),\n    @(\"64+26=90\", \"76-73=3\", \"73-12=61\", \"5+35=40\", \"43-15=28\"),\n    @(\"58-0=58\", \"24+48=72\", \"50+8=58\", \"48+38=86\", \"72-8=64\"),\n    @(\"16+2=18\", \"39+33=72\", \"98-93=5\", \"51-23=28\", \"81+5=86\"),\n    @(\"30+67=97\", \"55+32=87\", \"74-62=12\", \"65+30=95\", \"10+49=59\"),\n    @(\"20+61=81\", \"30+8=38\", \"9+22=31\", \"53+24=77\", \"8+54=62\"),\n    @(\"29+19=48\", \"83-75=8\", \"91-78=13\", \"14+14=28\", \"33+44=77\"),\n    @(\"21+40=61\", \"21+24=45\", \"77+0=77\", \"91-43=48\", \"46-9=37\"),\n    @(\"89-66=23\", \"77+0=77\", \"97-9=88\", \"19+25=44\", \"3+24=27\"),\n    @(\"54-7=47\", \"58+8=66\", \"10+26=36\", \"84-50=34\", \"82-75=7\"),\n    @(\"25+46=71\", \"37+26=63\", \"91-6=85\", \"25+44=69\", \"73-29=44\"),\n    @(\"12+65=77\", \"45+37=82\", \"54+45=99\", \"80-44=36\", \"18-12=6\"),\n    @(\"63+36=99\", \"47+30=77\", \"34+39=73\", \"76-65=11\", \"88-75=13\"),\n    @(\"61-29=32\", \"83+4=87\", \"4+11=15\", \"69+30=99\", \"51+15=66\"),\n    @(\"19+60=79\", \"61-43=18\", \"93-93=0\", \"54+26=80\", \"64-37=27\"),\n    @(\"2+19=21\", \"87-8=79\", \"9+61=70\", \"42+31=73\", \"80-55=25\"),\n    @(\"99-6=93\", \"59+9=68\", \"34+59=93\", \"94-32=62\", \"78-4=74\")\n)\n\nfor ($r = 1; $r -le $values.Count; $r++) {\n    $row = $values[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
